$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.927.37"
$ws.Range("E2").Value = "  -1.07%  "
$ws.Range("D3").Value = "3.084.95"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'521.41"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").Value = "'136.57"
$ws.Range("E6").Value = "  -2.68%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.087.36"
$ws.Range("E8").Value = "  +0.42%  "
$ws.Range("D9").Value = "'0.451"
$ws.Range("E9").Value = "  +2.63%  "
$ws.Range("D10").Value = "'7.35"
$ws.Range("E10").Value = "  +2.69%  "
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("D12").Value = "'0.398"
$ws.Range("E12").Value = "  +2.44%  "
$ws.Range("D13").Value = "3.619.32"
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("E14").Value = "  +1.52%  "
$ws.Range("D15").Value = "'25.37"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("D16").Value = "'0.0000161"
$ws.Range("E16").Value = "  -1.11%  "
$ws.Range("D17").Value = "57.100.96"
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("D18").Value = "3.083.30"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").Value = "'5.88"
$ws.Range("E19").Value = "  -2.98%  "
$ws.Range("D20").Value = "'12.46"
$ws.Range("E20").Value = "  -1.24%  "
$ws.Range("D21").Value = "'7.86"
$ws.Range("E21").Value = "  -0.74%  "
$ws.Range("D22").Value = "'347.23"
$ws.Range("E22").Value = "  +2.41%  "
$ws.Range("D23").Value = "'5.79"
$ws.Range("E23").Value = "  +1.51%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "'68.21"
$ws.Range("E25").Value = "  +2.04%  "
$ws.Range("D26").Value = "'0.497"
$ws.Range("E26").Value = "  -2.11%  "
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("D28").Value = "'1.01"
$ws.Range("E28").Value = "  +0.77%  "
$ws.Range("D29").Value = "0.0₃0880"
$ws.Range("E29").Value = "  -2.98%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").Value = "'7.28"
$ws.Range("E31").Value = "  +1.05%  "
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("D33").Value = "'5.90"
$ws.Range("E33").Value = "  -7.16%  "
$ws.Range("D34").Value = "'20.76"
$ws.Range("E34").Value = "  -0.41%  "
$ws.Range("D35").Value = "'4.91"
$ws.Range("E35").Value = "  +7.19%  "
$ws.Range("E36").Value = "  -3.07%  "
$ws.Range("D37").Value = "'159.33"
$ws.Range("E37").Value = "  +0.27%  "
$ws.Range("D38").Value = "'6.01"
$ws.Range("E38").Value = "  -1.57%  "
$ws.Range("D39").Value = "'25.80"
$ws.Range("E39").Value = "  -0.46%  "
$ws.Range("D40").Value = "'1.23"
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("D41").Value = "'0.0652"
$ws.Range("E41").Value = "  -1.61%  "
$ws.Range("D42").Value = "'1.58"
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").Value = "'4.02"
$ws.Range("E43").Value = "  +1.24%  "
$ws.Range("D44").Value = "'0.692"
$ws.Range("E44").Value = "  +1.82%  "
$ws.Range("D45").Value = "2.386.15"
$ws.Range("E45").Value = "  +4.75%  "
$ws.Range("E46").Value = "  -0.72%  "
$ws.Range("E47").Value = "  +0.15%  "
$ws.Range("D48").Value = "3.126.84"
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("D49").Value = "'0.0262"
$ws.Range("E49").Value = "  +0.70%  "
$ws.Range("D50").Value = "'0.959"
$ws.Range("E50").Value = "  -3.15%  "
$ws.Range("D51").Value = "'5.94"
$ws.Range("E51").Value = "  -1.84%  "
